$wb = $excel.ActiveWorkbook

# Locate the "I2C Cable" sheet so we can insert the new sheet right after it
$i2c = $wb.Worksheets.Item("I2C Cable")
$canBus = $wb.Worksheets.Item("CAN BUS Cable")

# Insert a new worksheet before "CAN BUS Cable" (i.e. right after "I2C Cable")
$newSheet = $wb.Worksheets.Add($canBus)
$newSheet.Name = "USB Cable"

$ws = $newSheet

$ws.Range("A1").Value = "For USB communication"

$ws.Range("A2").Value = "Manufacturer No."
$ws.Range("B2").Value = "Manufacturer"
$ws.Range("C2").Value = "Supplier"
$ws.Range("D2").Value = "Description"
$ws.Range("E2").Value = "Quantity"
$ws.Range("F2").Value = "Price"
$ws.Range("G2").Value = "Total"
$ws.Range("H2").Value = "note"

$ws.Range("A3").Value = "MUSBR-AHD2-241SK"
$ws.Range("B3").Value = "Amphenol"
$ws.Range("C3").Value = "Mouser"
$ws.Range("D3").Value = "USB Connectors IP67 Boot Style Hood Black for plug side"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 12.82
$ws.Range("G3").Formula = "=F3*E3"

$ws.Range("B4").Value = "Amphenol"
$ws.Range("C4").Value = "Mouser"
$ws.Range("D4").Value = "USB Connectors Rugged USB 3.0 A Vert Standard Shell"
$ws.Range("A4").Value = "MUSBR-3593-M0"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 15.07
$ws.Range("G4").Formula = "=F4*E4"
